# "Popravljen OdabirRestorana pdf, detaljnije uradjen glavni uml"
#
# On Sheet2 ("Odabir restorana" use case), the "Vezani zahtjevi" cell (B4)
# used to hold a long confirmation-flow description; that text was dropped
# (folded into the more detailed UML now), so B4 becomes "/" like the other
# not-applicable cells on this sheet. Row 4 no longer needs the tall wrapped
# layout, so its custom height shrinks back to the sheet default look.
# The view was also left scrolled/selected near the bottom of the table;
# reset it near the top, matching where editing actually continued.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# --- Cell content correction -------------------------------------------------
$ws.Range("B4").Value = "/"

# --- Row height: B4 no longer needs the tall wrapped row --------------------
$ws.Rows.Item(4).RowHeight = 18.75

# --- View/selection bookkeeping ----------------------------------------------
$ws.Range("D14").Select()
